# 11-7-2022 HAFIZ MESTRE.xlsx — sync with "deposito" inventory copy.
#
# The RIF_LIMPIO column (V) had been auto-incremented per row
# (...CI11, CI12, CI13, ...). The correct value is the same RIF for every
# line (...CI11), and the accidental CO_SENCAMER column (W, a stray "1" on
# every row) needs to be wiped out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("V" + $r).Value2 = "0005700003577CI11"
    $ws.Range("W" + $r).ClearContents()
}

# The sheet was scrolled/selected around C36:D36 while reviewing the bottom
# of the table, leaving an extra (empty) formatted row behind.
$ws.Rows.Item(36).RowHeight = 22.5
$ws.Range("C36:D36").Select()
